# TC06_CDS_Filter_ExprStrtgies-GWA.xlsx - Experimental Strategy and Study Data types - 13 Test cases
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New Cypher query text for each tab -----------------------------------

$filesTabQuery = "MATCH (f:file)-->(s:study)`nOPTIONAL MATCH (samp:sample)<--(f)`nOPTIONAL MATCH (samp)-->(p:participant)`nOPTIONAL MATCH (f)<--(g:genomic_info)`nOPTIONAL MATCH (p)<--(diag:diagnosis)`nWITH s, p, samp, f, g, diag, apoc.coll.flatten(COLLECT (apoc.text.split(f.experimental_strategy_and_data_subtypes,`"[;,]\s{0,1}`")), true) as es`nWHERE `"GWA`" IN es`nWITH DISTINCT f, s, p, samp`nRETURN`n    coalesce(f.file_name, '') as ``File Name``,`n    coalesce(s.study_name,'') as ``Study Name``,`n    coalesce(s.phs_accession,'') as ``Accession``,`n    coalesce(p.participant_id, '') as ``Participant ID``,`n    coalesce(samp.sample_id, '') as ``Sample ID``,`n    coalesce(f.file_type, '') as ``File Type```nORDER BY f.file_name LIMIT 100"

$samplesTabQuery = "MATCH (samp:sample)-->(p:participant)-->(s:study)`nOPTIONAL MATCH (samp)<--(f:file)`nOPTIONAL MATCH (f)<--(g:genomic_info)`nOPTIONAL MATCH (p)<--(diag:diagnosis)`nWITH s, p, samp, f, g, diag, apoc.coll.flatten(COLLECT (apoc.text.split(f.experimental_strategy_and_data_subtypes,`"[;,]\s{0,1}`")), true) as es`nWHERE `"GWA`" IN es`nWITH DISTINCT s, p, samp`nRETURN`n    coalesce(samp.sample_id, '') as ``Sample ID``,`n    coalesce(p.participant_id,'') as ``Participant ID``,`n    coalesce(s.study_name, '') as ``Study Name``,`n    coalesce(s.phs_accession,'') as ``Accession``,`n    coalesce(samp.sample_tumor_status,'') as ``Tumor``,`n    coalesce(samp.sample_type,'') as ``Analyte Type```nORDER BY samp.sample_id LIMIT 100"

$participantsTabQuery = "MATCH (p:participant)-->(s:study)`nOPTIONAL MATCH (samp:sample)-->(p)`nOPTIONAL MATCH (samp)<--(f:file)`nWITH p, samp, apoc.coll.flatten(COLLECT (apoc.text.split(f.experimental_strategy_and_data_subtypes,`"[;,]\s{0,1}`")), true) as es`nWHERE `"GWA`" IN es`nWITH p`nOPTIONAL MATCH (p)-->(s:study)`nOPTIONAL MATCH (samp:sample)-->(p)`nWITH s, p, apoc.coll.sort(collect(distinct coalesce(samp.sample_id, `"Not specified in data`"))) as samp`nRETURN `ncoalesce(p.participant_id,'') as ``Participant ID``,`ncoalesce(s.study_name, '') as ``Study Name``,`ncoalesce(s.phs_accession,'') as ``Accession``,`ncoalesce(p.gender,'') as ``Gender``,`ncoalesce(apoc.text.join(samp, ','), '') as ``Samples```nORDER BY p.participant_id LIMIT 100"

$statQuery = "CALL{`n    MATCH (p:participant)-->(s:study)`n    OPTIONAL MATCH (samp:sample)-->(p)`n    OPTIONAL MATCH (samp)<--(f:file)`n    OPTIONAL MATCH (f)<--(g:genomic_info)`n    OPTIONAL MATCH (p)<--(diag:diagnosis)`n    WITH s, p, samp, f, g, diag, apoc.coll.flatten(COLLECT (apoc.text.split(f.experimental_strategy_and_data_subtypes,`"[;,]\s{0,1}`")), true) as es`n    WHERE `"GWA`" IN es`n    RETURN `n        count(distinct p) AS num_participants`n}`nWITH num_participants`nCALL {`n    MATCH (samp:sample)-->(p:participant)-->(s:study)`n    OPTIONAL MATCH (samp)<--(f:file)`n    OPTIONAL MATCH (p)<--(diag:diagnosis)`n    OPTIONAL MATCH (f)<--(g:genomic_info)`n    OPTIONAL MATCH (p)<--(diag:diagnosis)`n    WITH s, p, samp, f, g, diag, apoc.coll.flatten(COLLECT (apoc.text.split(f.experimental_strategy_and_data_subtypes,`"[;,]\s{0,1}`")), true) as es`n    WHERE `"GWA`" IN es`n    RETURN `n        count(distinct samp) AS num_samples`n}`nWITH num_participants, num_samples`nCALL {`n    MATCH (f:file)-->(s:study)`n    OPTIONAL MATCH (f)<--(g:genomic_info)`n    OPTIONAL MATCH (samp:sample)<--(f)`n    OPTIONAL MATCH (p:participant)<--(samp)`n    OPTIONAL MATCH (p)<--(diag:diagnosis)`n    WITH s, p, samp, f, g, diag, apoc.coll.flatten(COLLECT (apoc.text.split(f.experimental_strategy_and_data_subtypes,`"[;,]\s{0,1}`")), true) as es`n    WHERE `"GWA`" IN es`n    RETURN `n        count(distinct s) AS num_studies,`n        count(distinct f) AS num_files`n}`nRETURN `n    num_studies AS Studies,`n    num_participants AS Participants,`n    num_samples AS Samples,`n    num_files AS ``Files``"

# --- Write the updated query text into the grid ----------------------------
# Row 2 = ParticipantsTab, Row 3 = SamplesTab, Row 4 = FilesTab
$ws.Range("B2").Value = $participantsTabQuery
$ws.Range("C2").Value = $statQuery

$ws.Range("B3").Value = $samplesTabQuery
$ws.Range("C3").Value = $statQuery

$ws.Range("B4").Value = $filesTabQuery
$ws.Range("C4").Value = $statQuery

# --- Formatting: bump base font size from 12pt to 14pt everywhere ---------
$ws.Range("A1:E4").Font.Size = 14
$ws.Range("B5:C5").Font.Size = 14
$ws.Range("C6").Font.Size = 14

# Ensure wrap text stays only on the query columns (B2:C4) plus the
# trailing helper cells (B5,C5,C6); everything else keeps no wrap.
$ws.Range("B2:C4").WrapText = $true
$ws.Range("B5:C5").WrapText = $true
$ws.Range("C6").WrapText = $true

# --- Row heights -------------------------------------------------------
$ws.Rows(2).RowHeight = 312
$ws.Rows(3).RowHeight = 388.5
$ws.Rows(4).RowHeight = 384.75

# --- View: scroll down so row 3 is at the top, select D4 ------------------
$ws.Range("D4").Select()
